$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "catalogos completos y pocas validaciones":
# Mark the "Catalogo Tipos Usuario" (row 17) and "Catalogo Dias Feriados" (row 18)
# sections as complete with an "OK" flag in column C, bold/size-14 text on a
# green fill.
$c17 = $ws.Range("C17")
$c17.Value = "OK"
$c17.Font.Bold = $true
$c17.Font.Size = 14
$c17.Interior.Color = 5296274

# Reuse the exact same formatting for C18 by copying C17's format instead of
# re-applying each property (keeps a single shared cell style).
$null = $c17.Copy()
$c18 = $ws.Range("C18")
$null = $c18.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$c18.Value = "OK"
$excel.CutCopyMode = $false

# Page is now set up for printing in portrait orientation.
$ws.PageSetup.Orientation = [Microsoft.Office.Interop.Excel.XlPageOrientation]::xlPortrait

# Leave the selection where the editor left it.
$null = $ws.Range("B20").Select()
